$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (GLD)
$ws.Range("D2").Value = 387.24
$ws.Range("E2").Value = 58.6
$ws.Range("F2").Value = 1.85
$ws.Range("H2").Value = 66
$ws.Range("K2").Value = 67
$ws.Range("N2").Value = 66.04328690552585

# Row 3 (NEM)
$ws.Range("D3").Value = 90.48
$ws.Range("E3").Value = 51.4
$ws.Range("F3").Value = 4.87
$ws.Range("I3").Value = 73
$ws.Range("J3").Value = 80
$ws.Range("K3").Value = 67
$ws.Range("N3").Value = 66.04328690552585

# Row 4 (GC=F)
$ws.Range("D4").Value = 4240
$ws.Range("E4").Value = 71.7
$ws.Range("F4").Value = 4.4
$ws.Range("J4").Value = 70
$ws.Range("K4").Value = 57
$ws.Range("N4").Value = 66.04328690552585
